$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 12, pushing existing rows 12-20 down to 13-21.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new weekly record.
$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(12, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(12, 3).Value = "Los Lagos"
$ws.Cells.Item(12, 4).Value = 44806
$ws.Cells.Item(12, 4).NumberFormat = $ws.Cells.Item(13, 4).NumberFormat
$ws.Cells.Item(12, 5).Value = 10
$ws.Cells.Item(12, 6).Value = 100112035
$ws.Cells.Item(12, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 70
$ws.Cells.Item(12, 11).Value = 23000
$ws.Cells.Item(12, 12).Value = 23000
$ws.Cells.Item(12, 13).Value = 23000
$ws.Cells.Item(12, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(12, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(12, 16).Value = 1533
$ws.Cells.Item(12, 17).Value = 15
$ws.Cells.Item(12, 18).Value = "Hortaliza"
